$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.158.02'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.08%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.902.86'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.57%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.10%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.49%  '

# Row 6
$ws.Range('E6').Value = '  +0.12%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5229'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.53%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3759'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.68%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07252'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.36%  '

# Row 10
$ws.Range('E10').Value = '  -0.57%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9021'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.59%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08493'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +11.44%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.918.69'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.47%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '95.07'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.22%  '

# Row 15
$ws.Range('E15').Value = '  +0.16%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.0000'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.07%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008636'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.32%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.53'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.44%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9999'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.14%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.195.06'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.16%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.066'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.26%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.147.39'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.18%  '

# Row 23
$ws.Range('E23').Value = '  +0.16%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.426'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.11%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.46'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.55%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.283'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.34%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.752'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.12%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.20'
$ws.Range('D28').Style = 'Normal'

# Row 29
$ws.Range('E29').Value = '  +0.26%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.816'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.99%  '

# Row 31
$ws.Range('E31').Value = '  -1.80%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09253'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.46%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8077'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.63%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05054'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.54%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.235'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.50%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.445'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.71%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.948'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.90%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.622'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.08%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5697'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.08%  '

# Row 40
$ws.Range('E40').Value = '  -0.27%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.075'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.12%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.024'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.35%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.632'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.23%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '116.52'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.36%  '

# Row 45
$ws.Range('E45').Value = '  +0.01%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4856'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.66%  '

# Row 47
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.16'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.75%  '

# Row 48
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.0000'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.15%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.616'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.12%  '

# Row 50
$ws.Range('E50').Value = '  +0.14%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '63.97'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.64%  '
